$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Cost / Stat table (K20:M47) - "Number Bought" / "Stat" / "Cost"
# -----------------------------------------------------------------
$ws.Range("K20").Value = "Number Bought"
$ws.Range("L20").Value = "Stat"
$ws.Range("M20").Value = "Cost"

$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 10
$ws.Range("M21").Value = 1000

$ws.Range("K22").Value = 1
$ws.Range("L22").Formula = "=L21+10"
$ws.Range("M22").Formula = "=M21*1.1"

for ($row = 23; $row -le 47; $row++) {
    $prev = $row - 1
    $ws.Range("K$row").Value = $row - 21
    $ws.Range("L$row").Formula = "=L$prev+10"
    $ws.Range("M$row").Formula = "=M$prev*1.1"
}

# -----------------------------------------------------------------
# New column widths for K, L, M
# -----------------------------------------------------------------
$ws.Columns.Item(11).ColumnWidth = 15.140625
$ws.Columns.Item(12).ColumnWidth = 7.85546875
$ws.Columns.Item(13).ColumnWidth = 15.7109375

# -----------------------------------------------------------------
# Thin border separating the new table (J column gets a left
# border, row 18 gets a divider row)
# -----------------------------------------------------------------
$ws.Columns.Item(10).Borders.Item(7).LineStyle = 1

$ws.Range("A18:D18").HorizontalAlignment = -4131
$ws.Range("C18").HorizontalAlignment = -4152
$ws.Range("D18").HorizontalAlignment = -4131

$ws.Range("J18").Borders.Item(7).LineStyle = 1
$ws.Range("J18").Borders.Item(8).LineStyle = 1

# -----------------------------------------------------------------
# Sheet view: scroll position + selection
# -----------------------------------------------------------------
$ws.Range("L40").Select()
$excel.ActiveWindow.Zoom = 100
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 5

Write-Output "done"
